# Updates cryptocurrency price/volume(1h) figures in the worksheet
# to reflect the latest data refresh (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.444.75"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "2.379.66"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'318.70"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "'109.25"
$ws.Range("E6").Value = "  -3.38%  "
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").Value = "'41.11"
$ws.Range("E10").Value = "  -3.70%  "
$ws.Range("D12").Value = "'8.58"
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("D14").Value = "'0.990"
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("D15").Value = "2.741.49"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "'15.52"
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("D17").Value = "2.374.52"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "45.386.86"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "'15.69"
$ws.Range("E19").Value = "  +16.90%  "
$ws.Range("D20").Value = "'7.34"
$ws.Range("E20").Value = "  -3.95%  "
$ws.Range("E21").Value = "  -1.07%  "
$ws.Range("D22").Value = "'3.73"
$ws.Range("E22").Value = "  +4.85%  "
$ws.Range("D23").Value = "'73.34"
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("D24").Value = "'261.48"
$ws.Range("E24").Value = "  -3.18%  "
$ws.Range("D25").Value = "'2.35"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'7.61"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("E29").Value = "  -4.54%  "
$ws.Range("D30").Value = "'22.48"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").Value = "'37.61"
$ws.Range("E32").Value = "  -4.38%  "
$ws.Range("D33").Value = "'167.86"
$ws.Range("E33").Value = "  -2.06%  "
$ws.Range("D34").Value = "'2.90"
$ws.Range("E34").Value = "  -1.69%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  -1.87%  "
$ws.Range("E37").Value = "  -3.69%  "
$ws.Range("E38").Value = "  +14.58%  "
$ws.Range("D39").Value = "'4.03"
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").Value = "'2.99"
$ws.Range("E40").Value = "  -2.82%  "
$ws.Range("E41").Value = "  -2.24%  "
$ws.Range("D42").Value = "'97.45"
$ws.Range("E42").Value = "  -7.44%  "
$ws.Range("E43").Value = "  -1.31%  "
$ws.Range("E44").Value = "  -4.02%  "
$ws.Range("D45").Value = "'12.95"
$ws.Range("E45").Value = "  -2.21%  "
$ws.Range("D46").Value = "1.854.87"
$ws.Range("E46").Value = "  +12.92%  "
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").Value = "'5.98"
$ws.Range("E48").Value = "  +3.28%  "
$ws.Range("D49").Value = "'83.53"
$ws.Range("E49").Value = "  +5.53%  "
$ws.Range("D50").Value = "'112.91"
$ws.Range("E50").Value = "  -3.24%  "
$ws.Range("E51").Value = "  -0.64%  "
